$d = $word.ActiveDocument

# Locate the final paragraph in the document body, which currently reads:
#   "A vertex is discovered the first time it is encountered during the
#    search, at which time it becomes nonwhite. "
$p = $d.Paragraphs.Last
$r = $p.Range

# Build the replacement paragraph content: the existing sentence stays as
# its own run, followed by new runs for "Gray and black " and a
# spell-check-flagged "verticies" (intentionally misspelled, matching the
# author's commit) wrapped in proofErr spellStart/spellEnd markers, and a
# trailing run containing a single space - mirroring the target OOXML.
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">A vertex is discovered the first time it is encountered during the search, at which time it becomes nonwhite. </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">Gray and black </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cstheme="minorHAnsi"/></w:rPr><w:t>verticies</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# InsertXML replaces the paragraph's text content but leaves the original
# paragraph mark (and its pPr/paraId/rsid attributes) behind as a new,
# now-empty, trailing paragraph. Insert first (discard the returned value
# so it doesn't get echoed to the output stream)...
$null = $r.InsertXML($xml)

# ...then merge that freshly-inserted paragraph back into the one that
# kept the original paragraph mark, by deleting the paragraph mark that
# sits between them. This restores the original paragraph count while
# preserving the original paragraph's identity (pStyle/numPr/paraId/rsids).
$newPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$mark = $d.Range($newPara.Range.End - 1, $newPara.Range.End)
$mark.Delete()
